$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 13747.5
$ws.Range("I9").Value = 1080
$ws.Range("J9").Value = 51750
$ws.Range("K9").Value = 1080
$ws.Range("L9").Value = 51750
$ws.Range("M9").Value = -911
$ws.Range("N9").Value = -52088

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5732.0586
$ws.Range("I19").Value = 642.25
$ws.Range("J19").Value = 7298.154
$ws.Range("K19").Value = 642.25
$ws.Range("L19").Value = 7298.154
$ws.Range("M19").Value = -467.25
$ws.Range("N19").Value = -7648.154

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 168207.86
$ws.Range("I92").Value = 95412.7
$ws.Range("K92").Value = 95412.7
$ws.Range("M92").Value = -94164.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 21663
$ws.Range("I106").Value = 21663
$ws.Range("K106").Value = 21663
$ws.Range("M106").Value = -21032

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1682.0189
$ws.Range("I132").Value = 1514.2094
$ws.Range("J132").Value = 2403.6
$ws.Range("K132").Value = 4542.6282
$ws.Range("L132").Value = 7210.799999999999
$ws.Range("M132").Value = -2012.6282
$ws.Range("N132").Value = -12270.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1089.5
$ws.Range("I135").Value = 1068.871
$ws.Range("J135").Value = 1160.5555
$ws.Range("K135").Value = 9619.839
$ws.Range("L135").Value = 10444.9995
$ws.Range("M135").Value = -7084.839
$ws.Range("N135").Value = -15514.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1827.3334
$ws.Range("I137").Value = 1191.2572
$ws.Range("K137").Value = 3573.7716
$ws.Range("M137").Value = -1023.7716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2386.2295
$ws.Range("I138").Value = 802.3871
$ws.Range("K138").Value = 2407.1613
$ws.Range("M138").Value = 2732.8387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9687.754999999999
$ws.Range("I32").Value = 4986.1777
$ws.Range("K32").Value = 4986.1777
$ws.Range("M32").Value = -4699.1777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 23760.334
$ws.Range("I36").Value = 10626
$ws.Range("K36").Value = 10626
$ws.Range("M36").Value = -10280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1361.4651
$ws.Range("I132").Value = 1464.1621
$ws.Range("K132").Value = 4392.4863
$ws.Range("M132").Value = -1862.4863

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2025.2174
$ws.Range("I20").Value = 1531.8572
$ws.Range("J20").Value = 2792.6667
$ws.Range("K20").Value = 1531.8572
$ws.Range("L20").Value = 2792.6667
$ws.Range("M20").Value = -1284.8572
$ws.Range("N20").Value = -3286.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 34497.5
$ws.Range("J21").Value = 34497.5
$ws.Range("L21").Value = 34497.5
$ws.Range("N21").Value = -34969.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 31332
$ws.Range("J95").Value = 31332
$ws.Range("L95").Value = 31332
$ws.Range("N95").Value = -36824

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2589.4285
$ws.Range("I105").Value = 2005.4
$ws.Range("K105").Value = 2005.4
$ws.Range("M105").Value = -258.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15022.5
$ws.Range("I107").Value = 16973.588
$ws.Range("J107").Value = 3966.3333
$ws.Range("K107").Value = 16973.588
$ws.Range("L107").Value = 3966.3333
$ws.Range("M107").Value = -15053.588
$ws.Range("N107").Value = -7806.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2072.25
$ws.Range("I16").Value = 748
$ws.Range("K16").Value = 748
$ws.Range("M16").Value = -461

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 19300
$ws.Range("J37").Value = 21200
$ws.Range("L37").Value = 21200
$ws.Range("N37").Value = -21414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 1521
$ws.Range("I38").Value = 1000
$ws.Range("J38").Value = 2042
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = 2042
$ws.Range("M38").Value = -623
$ws.Range("N38").Value = -2796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 1521
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2042
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2042
$ws.Range("M46").Value = -789
$ws.Range("N46").Value = -2464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 39617
$ws.Range("I107").Value = 70266.234
$ws.Range("J107").Value = 3395.182
$ws.Range("K107").Value = 70266.234
$ws.Range("L107").Value = 3395.182
$ws.Range("M107").Value = -68346.234
$ws.Range("N107").Value = -7235.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2072.25
$ws.Range("I113").Value = 748
$ws.Range("K113").Value = 748
$ws.Range("M113").Value = 1422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 58975
$ws.Range("J120").Value = 58975
$ws.Range("L120").Value = 58975
$ws.Range("N120").Value = -66233

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2765.5881
$ws.Range("I134").Value = 2058.7273
$ws.Range("K134").Value = 6176.1819
$ws.Range("M134").Value = -3641.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1953.5
$ws.Range("I17").Value = 1573.6666
$ws.Range("J17").Value = 2333.3333
$ws.Range("K17").Value = 4720.9998
$ws.Range("L17").Value = 6999.999899999999
$ws.Range("M17").Value = -4551.9998
$ws.Range("N17").Value = -7337.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1627.9524
$ws.Range("J107").Value = 1732.4736
$ws.Range("L107").Value = 5197.4208
$ws.Range("N107").Value = -9037.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4402.1665
$ws.Range("J113").Value = 4402.1665
$ws.Range("L113").Value = 13206.4995
$ws.Range("N113").Value = -17546.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 5830.1055
$ws.Range("I129").Value = 1333.3334
$ws.Range("J129").Value = 6673.25
$ws.Range("K129").Value = 4000.0002
$ws.Range("L129").Value = 20019.75
$ws.Range("M129").Value = 999.9998000000001
$ws.Range("N129").Value = -30019.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3734.8333
$ws.Range("J131").Value = 5988.4
$ws.Range("L131").Value = 17965.2
$ws.Range("N131").Value = -28045.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18693.422
$ws.Range("I102").Value = 23489.512
$ws.Range("K102").Value = 23489.512
$ws.Range("M102").Value = -21867.512

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6236.154
$ws.Range("I113").Value = 6338.7827
$ws.Range("K113").Value = 6338.7827
$ws.Range("M113").Value = -4168.7827

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 59170.5
$ws.Range("I122").Value = 92711.3
$ws.Range("J122").Value = 3269.1667
$ws.Range("K122").Value = 278133.9
$ws.Range("L122").Value = 9807.500100000001
$ws.Range("M122").Value = -275683.9
$ws.Range("N122").Value = -14707.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 50956.156
$ws.Range("I126").Value = 78532.25
$ws.Range("J126").Value = 3682.8572
$ws.Range("K126").Value = 235596.75
$ws.Range("L126").Value = 11048.5716
$ws.Range("M126").Value = -233126.75
$ws.Range("N126").Value = -15988.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2038.5883
$ws.Range("I46").Value = 776.3333
$ws.Range("J46").Value = 2309.0715
$ws.Range("K46").Value = 776.3333
$ws.Range("L46").Value = 2309.0715
$ws.Range("M46").Value = -588.3333
$ws.Range("N46").Value = -2685.0715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 16299.4
$ws.Range("J58").Value = 19124.5
$ws.Range("L58").Value = 19124.5
$ws.Range("N58").Value = -19644.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 140857.42
$ws.Range("J127").Value = 140857.42
$ws.Range("L127").Value = 140857.42
$ws.Range("N127").Value = -150777.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 23950
$ws.Range("I32").Value = 23950
$ws.Range("K32").Value = 23950
$ws.Range("M32").Value = -23633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 45282.43
$ws.Range("J95").Value = 45282.43
$ws.Range("L95").Value = 45282.43
$ws.Range("N95").Value = -50774.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1892.4783
$ws.Range("I126").Value = 1876.95
$ws.Range("K126").Value = 5630.85
$ws.Range("M126").Value = -3160.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 129476
$ws.Range("J140").Value = 129476
$ws.Range("L140").Value = 129476
$ws.Range("N140").Value = -139836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 58199.355
$ws.Range("J141").Value = 58199.355
$ws.Range("L141").Value = 58199.355
$ws.Range("N141").Value = -68559.35500000001
